# CaveStory schedule workbook - "Add files via upload" / 10.28 update
#
# Content changes on sheet "계획표":
#  - E7:  "대장간(무기 제작)\n무기장착"  -> "가구 구현"          (furniture feature replaces old smithy task)
#  - F7:  "주민 (애니메이션) & 영웅"      -> "정산창 & 종족 진화"  (swapped with H7)
#  - H7:  "정산창 & 종족 진화"           -> "주민 (애니메이션) & 영웅"  (swapped with F7)
#  - G14: (empty) -> "정산창 & 종족 진화"        (new "추후 작업"/future-work row entry, bordered like the table)
#  - H14: (empty) -> "주민 (애니메이션) & 영웅"  (new "추후 작업"/future-work row entry, bordered like the table)
# Plus the active selection moves to E4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 edits ---
$ws.Range("E7").Value = "가구 구현"
$ws.Range("F7").Value = "정산창 & 종족 진화"
$ws.Range("H7").Value = "주민 (애니메이션) & 영웅"

# --- New entries in the "추후 작업" (future work) row, matching the bordered/centered table style ---
$futureRange = $ws.Range("G14:H14")
$ws.Range("G14").Value = "정산창 & 종족 진화"
$ws.Range("H14").Value = "주민 (애니메이션) & 영웅"
$futureRange.Borders.LineStyle = 1
$futureRange.HorizontalAlignment = -4108
$futureRange.VerticalAlignment = -4108
$futureRange.WrapText = $false

# --- Selection moves to E4 ---
$ws.Range("E4").Select()
